# "Generate Report for Handoff":
#  - the source markdown file tracked by this status report is renamed
#    (new content hash) from cf39f490-...md to 48566a70-...md, and its
#    handoff package (.xlf) is re-generated with a new hash + later
#    timestamp;
#  - a second source file (ffff26b2e6c4-...md) is now also tracked, with
#    its own "Ready for handoff" row sitting right above the constant
#    ".localization-config" row (which shifts down one row on every
#    sheet).
#
# This touches all three worksheets (Overview, zh-cn, de-de): each one
# gains a new row 3, and its old row 3 (".localization-config") becomes
# row 4.

$wb = $excel.ActiveWorkbook

$oldMd  = "cf39f490-87b4-4cea-9542-1190327a7289.md"
$newMd  = "48566a70-0a28-4fce-8ad0-9368ac6f1432.md"
$newMd2 = "ffff26b2e6c4-cc09-49e5-abf2-c9e0f58efd83.md"
$cfgFile = ".localization-config"

$repoBase   = "https://github.com/OpenLocalizationTest/oltest/blob/d3f847c598ba64673e6ab520c4b85b066a8afa75"
$zhHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/24a91192cafd69542b4368a523ab1731ed41d601/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht"
$deHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b9bfbdf00e6156b441151a8230ac5882aeecafd0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht"

$newXlfZh = "48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.zh-cn.xlf"
$newXlfDe = "48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.de-de.xlf"

$zhHandoffTime = "2016-02-22 05:07:48"
$deHandoffTime = "2016-02-22 05:08:03"
$epoch = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# push the constant config row down from row 3 -> row 4
$wsOverview.Range("A4").Value = $cfgFile
$wsOverview.Range("B4").Value = "Not to be localized"
$wsOverview.Range("C4").Value = "Not to be localized"

# row 2: source file renamed
$wsOverview.Range("A2").Value = $newMd

# row 3: newly tracked source file
$wsOverview.Range("A3").Value = $newMd2
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# hyperlinks shift with the rows, so rebuild the whole set for this sheet
$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "$repoBase/e2e/$newMd", $null, $null, $newMd)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "$repoBase/e2e/$newMd2", $null, $null, $newMd2)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "$repoBase/$cfgFile", $null, $null, $cfgFile)

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# push the constant config row down from row 3 -> row 4
$wsZh.Range("A4").Value = $cfgFile
$wsZh.Range("B4").Value = "Not to be localized"
$wsZh.Range("C4").Value = ""
$wsZh.Range("D4").Value = $epoch
$wsZh.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("G4").Value = $epoch
$wsZh.Range("H4").Value = "Ignored"

# row 2: source file renamed + new handoff package/time
$wsZh.Range("A2").Value = $newMd
$wsZh.Range("B2").Value = "Ready for handoff"
$wsZh.Range("C2").Value = $newXlfZh
$wsZh.Range("D2").Value = $zhHandoffTime
$wsZh.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("G2").Value = $epoch
$wsZh.Range("H2").Value = "Include"

# row 3: newly tracked source file, handed off together with row 2's file
$wsZh.Range("A3").Value = $newMd2
$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("C3").Value = $newXlfZh
$wsZh.Range("D3").Value = $zhHandoffTime
$wsZh.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("G3").Value = $epoch
$wsZh.Range("H3").Value = "Include"

# hyperlinks shift with the rows, so rebuild the whole set for this sheet
$wsZh.Range("A1").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "$repoBase/e2e/$newMd", $null, $null, $newMd)
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "$zhHandoffBase/$newXlfZh", $null, $null, $newXlfZh)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "$repoBase/e2e/$newMd2", $null, $null, $newMd2)
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), "$zhHandoffBase/$newXlfZh", $null, $null, $newXlfZh)
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "$repoBase/$cfgFile", $null, $null, $cfgFile)

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# push the constant config row down from row 3 -> row 4
$wsDe.Range("A4").Value = $cfgFile
$wsDe.Range("B4").Value = "Not to be localized"
$wsDe.Range("C4").Value = ""
$wsDe.Range("D4").Value = $epoch
$wsDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("G4").Value = $epoch
$wsDe.Range("H4").Value = "Ignored"

# row 2: source file renamed + new handoff package/time
$wsDe.Range("A2").Value = $newMd
$wsDe.Range("B2").Value = "Ready for handoff"
$wsDe.Range("C2").Value = $newXlfDe
$wsDe.Range("D2").Value = $deHandoffTime
$wsDe.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("G2").Value = $epoch
$wsDe.Range("H2").Value = "Include"

# row 3: newly tracked source file, handed off together with row 2's file
$wsDe.Range("A3").Value = $newMd2
$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("C3").Value = $newXlfDe
$wsDe.Range("D3").Value = $deHandoffTime
$wsDe.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("G3").Value = $epoch
$wsDe.Range("H3").Value = "Include"

# hyperlinks shift with the rows, so rebuild the whole set for this sheet
$wsDe.Range("A1").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$repoBase/e2e/$newMd", $null, $null, $newMd)
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "$deHandoffBase/$newXlfDe", $null, $null, $newXlfDe)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "$repoBase/e2e/$newMd2", $null, $null, $newMd2)
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), "$deHandoffBase/$newXlfDe", $null, $null, $newXlfDe)
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "$repoBase/$cfgFile", $null, $null, $cfgFile)

Write-Output "Handoff report regenerated: added $newMd2, refreshed $newMd handoff package."
